# Update Leve profit calculations across multiple job sheets
# (currentAveragePrice / NQ / HQ price & profit columns H:N)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 153.25
$ws.Range("I5").Value = 166
$ws.Range("J5").Value = 115
$ws.Range("K5").Value = 166
$ws.Range("L5").Value = 115
$ws.Range("M5").Value = -51
$ws.Range("N5").Value = -345

# Row 9
$ws.Range("H9").Value = 8462.5
$ws.Range("I9").Value = 11242.777
$ws.Range("J9").Value = 121.666664
$ws.Range("K9").Value = 11242.777
$ws.Range("L9").Value = 121.666664
$ws.Range("M9").Value = -11073.777

# Row 12
$ws.Range("H12").Value = 949.5
$ws.Range("I12").Value = 299
$ws.Range("J12").Value = 1339.8
$ws.Range("K12").Value = 299
$ws.Range("L12").Value = 1339.8
$ws.Range("M12").Value = -129

# Row 62
$ws.Range("H62").Value = 2391.7778
$ws.Range("I62").Value = 2391.7778
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2391.7778
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1767.7778

# Row 64
$ws.Range("H64").Value = 7882.4287
$ws.Range("I64").Value = 7796.25
$ws.Range("J64").Value = 7997.3335
$ws.Range("K64").Value = 7796.25
$ws.Range("L64").Value = 7997.3335
$ws.Range("M64").Value = -7548.25

# Row 65
$ws.Range("H65").Value = 2391.7778
$ws.Range("I65").Value = 2391.7778
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11958.889
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -8838.888999999999

# Row 67
$ws.Range("H67").Value = 7882.4287
$ws.Range("I67").Value = 7796.25
$ws.Range("J67").Value = 7997.3335
$ws.Range("K67").Value = 7796.25
$ws.Range("L67").Value = 7997.3335
$ws.Range("M67").Value = -6938.25

# Row 70
$ws.Range("H70").Value = 3398.204
$ws.Range("I70").Value = 1264.2222
$ws.Range("J70").Value = 9307.691999999999
$ws.Range("K70").Value = 3792.6666
$ws.Range("L70").Value = 27923.076
$ws.Range("M70").Value = -3522.6666
$ws.Range("N70").Value = -28463.076

# Row 73
$ws.Range("H73").Value = 3398.204
$ws.Range("I73").Value = 1264.2222
$ws.Range("J73").Value = 9307.691999999999
$ws.Range("K73").Value = 3792.6666
$ws.Range("L73").Value = 27923.076
$ws.Range("M73").Value = -2856.6666
$ws.Range("N73").Value = -29795.076

# Row 129
$ws.Range("H129").Value = 998.7222
$ws.Range("I129").Value = 898.5333000000001
$ws.Range("J129").Value = 1499.6666
$ws.Range("K129").Value = 2695.5999
$ws.Range("L129").Value = 4498.9998
$ws.Range("M129").Value = 2304.4001
$ws.Range("N129").Value = -14498.9998

# Row 137
$ws.Range("H137").Value = 1674.6666
$ws.Range("I137").Value = 1426.05
$ws.Range("J137").Value = 2385
$ws.Range("K137").Value = 4278.15
$ws.Range("L137").Value = 7155
$ws.Range("M137").Value = -1728.15
$ws.Range("N137").Value = -12255

# Row 138
$ws.Range("H138").Value = 4255.7837
$ws.Range("I138").Value = 4313.1763
$ws.Range("J138").Value = 4207
$ws.Range("K138").Value = 12939.5289
$ws.Range("L138").Value = 12621
$ws.Range("M138").Value = -7799.528900000001
$ws.Range("N138").Value = -22901

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 865.58826
$ws.Range("I2").Value = 662.8
$ws.Range("J2").Value = 2386.5
$ws.Range("K2").Value = 662.8
$ws.Range("L2").Value = 2386.5
$ws.Range("M2").Value = -549.8

# Row 32
$ws.Range("H32").Value = 3674.1738
$ws.Range("I32").Value = 3159.3635
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 3159.3635
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -2872.3635

# Row 61
$ws.Range("H61").Value = 5174.838
$ws.Range("I61").Value = 2331.9666
$ws.Range("J61").Value = 17358.572
$ws.Range("K61").Value = 2331.9666
$ws.Range("L61").Value = 17358.572
$ws.Range("M61").Value = -2119.9666

# Row 88
$ws.Range("H88").Value = 6526.1577
$ws.Range("I88").Value = 1207.8
$ws.Range("J88").Value = 8425.571
$ws.Range("K88").Value = 1207.8
$ws.Range("L88").Value = 8425.571
$ws.Range("M88").Value = -801.8
$ws.Range("N88").Value = -9237.571

# Row 91
$ws.Range("H91").Value = 6526.1577
$ws.Range("I91").Value = 1207.8
$ws.Range("J91").Value = 8425.571
$ws.Range("K91").Value = 1207.8
$ws.Range("L91").Value = 8425.571
$ws.Range("M91").Value = 196.2
$ws.Range("N91").Value = -11233.571

# Row 116
$ws.Range("H116").Value = 865.58826
$ws.Range("I116").Value = 662.8
$ws.Range("J116").Value = 2386.5
$ws.Range("K116").Value = 662.8
$ws.Range("L116").Value = 2386.5
$ws.Range("M116").Value = 1631.2

# Row 132
$ws.Range("H132").Value = 1523.7142
$ws.Range("I132").Value = 1464.9474
$ws.Range("J132").Value = 2082
$ws.Range("K132").Value = 4394.8422
$ws.Range("L132").Value = 6246
$ws.Range("M132").Value = -1864.8422
$ws.Range("N132").Value = -11306

# Row 136
$ws.Range("H136").Value = 5174.838
$ws.Range("I136").Value = 2331.9666
$ws.Range("J136").Value = 17358.572
$ws.Range("K136").Value = 6995.899800000001
$ws.Range("L136").Value = 52075.716
$ws.Range("M136").Value = -4445.899800000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 865.58826
$ws.Range("I3").Value = 662.8
$ws.Range("J3").Value = 2386.5
$ws.Range("K3").Value = 662.8
$ws.Range("L3").Value = 2386.5
$ws.Range("M3").Value = -548.8

# Row 99
$ws.Range("H99").Value = 3893.5833
$ws.Range("I99").Value = 3201.7144
$ws.Range("J99").Value = 4862.2
$ws.Range("K99").Value = 3201.7144
$ws.Range("L99").Value = 4862.2
$ws.Range("M99").Value = -1703.7144
$ws.Range("N99").Value = -7858.2

# Row 134
$ws.Range("H134").Value = 6162.8335
$ws.Range("I134").Value = 5847.905
$ws.Range("J134").Value = 7265.0835
$ws.Range("K134").Value = 17543.715
$ws.Range("L134").Value = 21795.2505
$ws.Range("M134").Value = -15008.715

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 125.947365
$ws.Range("I7").Value = 121.57143
$ws.Range("J7").Value = 138.2
$ws.Range("K7").Value = 121.57143
$ws.Range("L7").Value = 138.2
$ws.Range("M7").Value = -8.571430000000007

# Row 22
$ws.Range("H22").Value = 725
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1150
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1150
$ws.Range("M22").Value = 50

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 473244.3
$ws.Range("I4").Value = 569056.4
$ws.Range("J4").Value = 169839.5
$ws.Range("K4").Value = 1707169.2
$ws.Range("L4").Value = 509518.5
$ws.Range("M4").Value = -1707057.2
$ws.Range("N4").Value = -509742.5

# Row 132
$ws.Range("H132").Value = 2154.25
$ws.Range("I132").Value = 2066.1667
$ws.Range("J132").Value = 2418.5
$ws.Range("K132").Value = 18595.5003
$ws.Range("L132").Value = 21766.5
$ws.Range("M132").Value = -16065.5003
$ws.Range("N132").Value = -26826.5

# Row 137
$ws.Range("H137").Value = 3996.5715
$ws.Range("I137").Value = 2812.3
$ws.Range("J137").Value = 6957.25
$ws.Range("K137").Value = 8436.900000000001
$ws.Range("L137").Value = 20871.75
$ws.Range("M137").Value = -3336.900000000001
$ws.Range("N137").Value = -31071.75

# Row 141
$ws.Range("H141").Value = 9035.5
$ws.Range("I141").Value = 7040.5713
$ws.Range("J141").Value = 23000
$ws.Range("K141").Value = 21121.7139
$ws.Range("L141").Value = 69000
$ws.Range("M141").Value = -15941.7139

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 1346.2
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1346.2
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1346.2
$ws.Range("N107").Value = -5186.2

# Row 122
$ws.Range("H122").Value = 2590.8518
$ws.Range("I122").Value = 2294.2222
$ws.Range("J122").Value = 3184.111
$ws.Range("K122").Value = 6882.6666
$ws.Range("L122").Value = 9552.332999999999
$ws.Range("M122").Value = -4432.6666
$ws.Range("N122").Value = -14452.333

# Row 132
$ws.Range("H132").Value = 2573.9
$ws.Range("I132").Value = 2637.6667
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7913.000100000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -5383.000100000001
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 7298.375
$ws.Range("I81").Value = 7298.375
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 14596.75
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -13535.75

# Row 84
$ws.Range("H84").Value = 7298.375
$ws.Range("I84").Value = 7298.375
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 72983.75
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -67679.75

# Row 113
$ws.Range("H113").Value = 626641.5600000001
$ws.Range("I113").Value = 1430412.9
$ws.Range("J113").Value = 1486.1111
$ws.Range("K113").Value = 4291238.699999999
$ws.Range("L113").Value = 4458.3333
$ws.Range("M113").Value = -4289068.699999999
$ws.Range("N113").Value = -8798.3333

# Row 132
$ws.Range("H132").Value = 4239.475
$ws.Range("I132").Value = 4193.788
$ws.Range("J132").Value = 4454.857
$ws.Range("K132").Value = 12581.364
$ws.Range("L132").Value = 13364.571
$ws.Range("M132").Value = -10051.364

# Row 133
$ws.Range("H133").Value = 80000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 80000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
